# fix check for current holdings
# Sets SSG-Bibliothek (U) and Laufend (V) to "DE-7" / "x" for the listed
# rows, and also sets Unikal (W) to "x" for rows where it changed in the
# source diff (all affected rows except 10 and 12, where W stays blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where U -> "DE-7" and V -> "x"
$rowsUV = @(9, 10, 12, 18, 20, 21, 25, 56, 92, 101, 105, 122)
foreach ($r in $rowsUV) {
    $ws.Range("U$r").Value = "DE-7"
    $ws.Range("V$r").Value = "x"
}

# Rows where W also -> "x"
$rowsW = @(9, 18, 20, 21, 25, 56, 92, 101, 105, 122)
foreach ($r in $rowsW) {
    $ws.Range("W$r").Value = "x"
}
